$wb = $excel.ActiveWorkbook

# New B (Cutoff) and C (Reaction_number) values for data rows 2-16 (A column values 0-14)
# Sheet "NBR" (first worksheet)
$nbrB = @(5,6,7,8,9,10,11,12,13,14,15,16,17,18,19)
$nbrC = @(801,798,795,795,795,796,796,796,789,788,542,541,534,533,527)

# Sheet "BAR" (second worksheet)
$barB = @(5,6,7,8,9,10,11,12,13,14,15,16,17,18,19)
$barC = @(707,704,704,704,704,699,698,697,705,706,704,702,699,698,697)

$wsNBR = $wb.Worksheets.Item("NBR")
$wsBAR = $wb.Worksheets.Item("BAR")

for ($i = 0; $i -lt 15; $i++) {
    $r = $i + 2
    $wsNBR.Cells.Item($r, 2).Value = $nbrB[$i]
    $wsNBR.Cells.Item($r, 3).Value = $nbrC[$i]

    $wsBAR.Cells.Item($r, 2).Value = $barB[$i]
    $wsBAR.Cells.Item($r, 3).Value = $barC[$i]
}

# Remove the now-obsolete trailing rows (17-20), shrinking the used range to A1:C16
$wsNBR.Rows("17:20").Delete()
$wsBAR.Rows("17:20").Delete()
